# Refresh cryptos list values (price + 1h volume change) per upstream source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.986.09'
$ws.Range('E2').Value = '  -2.33%  '
$ws.Range('D3').Value = '2.099.56'
$ws.Range('E3').Value = '  -1.09%  '
$ws.Range('E4').Value = '  -0.86%  '
$ws.Range('D5').Value = '''346.62'
$ws.Range('E5').Value = '  +2.50%  '
$ws.Range('E6').Value = '  -0.80%  '
$ws.Range('D7').Value = '''0.5159'
$ws.Range('E7').Value = '  -1.86%  '
$ws.Range('D8').Value = '''0.4424'
$ws.Range('D9').Value = '''0.09383'
$ws.Range('E9').Value = '  +2.77%  '
$ws.Range('D10').Value = '''52.14'
$ws.Range('E10').Value = '  -4.80%  '
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('D12').Value = '''25.23'
$ws.Range('E12').Value = '  +2.86%  '
$ws.Range('D13').Value = '2.102.24'
$ws.Range('E13').Value = '  -0.83%  '
$ws.Range('D14').Value = '''6.749'
$ws.Range('E14').Value = '  -1.68%  '
$ws.Range('D15').Value = '''8.169'
$ws.Range('E15').Value = '  +0.08%  '
$ws.Range('D16').Value = '''99.51'
$ws.Range('E16').Value = '  +2.25%  '
$ws.Range('D17').Value = '''0.00001172'
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('D19').Value = '''20.59'
$ws.Range('E19').Value = '  +5.61%  '
$ws.Range('D20').Value = '''0.06681'
$ws.Range('E20').Value = '  -0.23%  '
$ws.Range('E21').Value = '  -0.82%  '
$ws.Range('D22').Value = '''6.217'
$ws.Range('E22').Value = '  -1.72%  '
$ws.Range('D23').Value = '30.088.27'
$ws.Range('E23').Value = '  -2.15%  '
$ws.Range('D24').Value = '''12.62'
$ws.Range('E24').Value = '  -2.35%  '
$ws.Range('D25').Value = '''2.331'
$ws.Range('E25').Value = '  -1.22%  '
$ws.Range('D26').Value = '2.347.27'
$ws.Range('E26').Value = '  -0.88%  '
$ws.Range('D27').Value = '''21.95'
$ws.Range('E27').Value = '  -2.15%  '
$ws.Range('D28').Value = '''2.557'
$ws.Range('E28').Value = '  -0.53%  '
$ws.Range('D29').Value = '''162.09'
$ws.Range('E29').Value = '  -1.69%  '
$ws.Range('D30').Value = '''133.48'
$ws.Range('E30').Value = '  -1.01%  '
$ws.Range('D31').Value = '''1.171'
$ws.Range('E31').Value = '  -3.40%  '
$ws.Range('D32').Value = '''0.1061'
$ws.Range('E32').Value = '  -1.37%  '
$ws.Range('D33').Value = '''1.643'
$ws.Range('E33').Value = '  -1.94%  '
$ws.Range('D34').Value = '''6.225'
$ws.Range('E34').Value = '  -2.49%  '
$ws.Range('D35').Value = '''3.960'
$ws.Range('E35').Value = '  +0.42%  '
$ws.Range('D36').Value = '''6.203'
$ws.Range('E36').Value = '  +5.35%  '
$ws.Range('D37').Value = '''10.08'
$ws.Range('E37').Value = '  -5.66%  '
$ws.Range('E38').Value = '  -3.14%  '
$ws.Range('D39').Value = '''0.06781'
$ws.Range('E39').Value = '  -1.53%  '
$ws.Range('D40').Value = '''0.2276'
$ws.Range('E40').Value = '  -2.41%  '
$ws.Range('D41').Value = '''0.6951'
$ws.Range('E41').Value = '  +0.34%  '
$ws.Range('D42').Value = '''12.51'
$ws.Range('E42').Value = '  -1.33%  '
$ws.Range('D43').Value = '''1.311'
$ws.Range('E43').Value = '  +3.88%  '
$ws.Range('D44').Value = '''0.6621'
$ws.Range('E44').Value = '  +1.78%  '
$ws.Range('D45').Value = '''14.16'
$ws.Range('E45').Value = '  -6.26%  '
$ws.Range('D46').Value = '''2.279'
$ws.Range('E46').Value = '  -1.75%  '
$ws.Range('E47').Value = '  -1.82%  '
$ws.Range('E48').Value = '  -4.89%  '
$ws.Range('E49').Value = '  -2.98%  '
$ws.Range('D50').Value = '''82.07'
$ws.Range('E50').Value = '  -1.84%  '
$ws.Range('D51').Value = '''0.07207'
$ws.Range('E51').Value = '  -1.34%  '
